$d = $word.ActiveDocument

# Move to the end of the document's last paragraph ("For this project, the
# only external module I plan on using is pygame") and append two new
# paragraphs after it, matching the existing body-text formatting
# (sz/szCs 24, line spacing 480 auto, rtl 0) which Word inherits
# automatically from the paragraph we split off of.

$tail = $d.Paragraphs.Last.Range
$tail.Collapse(0)
$tail.InsertParagraphAfter()

$p1 = $d.Paragraphs.Last.Range
$p1.Collapse(0)
$p1.InsertAfter("TP2 Update:")

$p1end = $d.Paragraphs.Last.Range
$p1end.Collapse(0)
$p1end.InsertParagraphAfter()

$p2 = $d.Paragraphs.Last.Range
$p2.Collapse(0)
$p2.InsertAfter("I am thinking about making the game not end until the player loses all their pikmin to enemies, but I am not sure about this feature so I will just include it here in case I end up implementing it.")
